$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.048.57"
$ws.Range("E2").Value = "  -0.56%  "

$ws.Range("D3").Value = "1.651.58"
$ws.Range("E3").Value = "  -0.45%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.50"
$ws.Range("E5").Value = "  +0.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5266"
$ws.Range("E6").Value = "  +2.03%  "

$ws.Range("E8").Value = "  -1.76%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06321"
$ws.Range("E9").Value = "  +0.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.38"
$ws.Range("E10").Value = "  -1.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07792"
$ws.Range("E11").Value = "  +0.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.502"
$ws.Range("E12").Value = "  +0.52%  "

$ws.Range("D13").Value = "1.651.95"
$ws.Range("E13").Value = "  +0.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5490"
$ws.Range("E14").Value = "  +0.66%  "

$ws.Range("D15").Value = "0.0₅8203"
$ws.Range("E15").Value = "  +1.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.42"
$ws.Range("E16").Value = "  +0.85%  "

$ws.Range("D17").Value = "26.064.55"
$ws.Range("E17").Value = "  -0.54%  "

$ws.Range("E18").Value = "  -0.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.579"
$ws.Range("E19").Value = "  -0.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.46"
$ws.Range("E20").Value = "  -0.32%  "

$ws.Range("E21").Value = "  -0.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.043"
$ws.Range("E22").Value = "  +0.52%  "

$ws.Range("E23").Value = "  -0.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "142.73"
$ws.Range("E24").Value = "  +2.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1239"
$ws.Range("E25").Value = "  +1.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.226"
$ws.Range("E26").Value = "  -0.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.06"
$ws.Range("E27").Value = "  -0.55%  "

$ws.Range("E28").Value = "  -0.82%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05827"
$ws.Range("E29").Value = "  -2.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.546"
$ws.Range("E31").Value = "  -0.32%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.264"
$ws.Range("E32").Value = "  +0.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.589"
$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.408"
$ws.Range("E34").Value = "  -0.56%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9439"
$ws.Range("E35").Value = "  -2.04%  "

$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.776"
$ws.Range("E36").Value = "  +0.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5742"
$ws.Range("E37").Value = "  +1.11%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01607"
$ws.Range("E38").Value = "  +0.96%  "

$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.746"
$ws.Range("E39").Value = "  -5.41%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8425"
$ws.Range("E40").Value = "  -1.74%  "

$ws.Range("E41").Value = "  -0.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "103.58"
$ws.Range("E42").Value = "  +3.19%  "

$ws.Range("D43").Value = "1.028.88"
$ws.Range("E43").Value = "  +1.75%  "

$ws.Range("D44").Value = "1.794.74"
$ws.Range("E44").Value = "  -0.33%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "56.99"
$ws.Range("E45").Value = "  +0.57%  "

$ws.Range("E46").Value = "  -0.45%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4321"
$ws.Range("E47").Value = "  +3.25%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.842"
$ws.Range("E48").Value = "  -2.49%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05144"
$ws.Range("E49").Value = "  -0.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.467"
$ws.Range("E50").Value = "  +1.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09634"
$ws.Range("E51").Value = "  -0.56%  "
